$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 1063; this shifts the existing
# rows 1063-1083 down to 1065-1085 and copies formatting (incl.
# the date style on column D) from the row above the insertion point.
$ws.Rows("1063:1064").Insert()

# Fill in the two new records (most recent observations) in the
# newly inserted rows 1063 and 1064.
$ws.Range("A1063").Value = 5
$ws.Range("B1063").Value = "Macroferia Regional de Talca"
$ws.Range("C1063").Value = "Maule"
$ws.Range("D1063").Value = 45239
$ws.Range("E1063").Value = 7
$ws.Range("F1063").Value = 100112004
$ws.Range("G1063").Value = "Cebolla"
$ws.Range("H1063").Value = "Sin especificar"
$ws.Range("I1063").Value = "1a nueva(o)"
$ws.Range("J1063").Value = 20000
$ws.Range("K1063").Value = 3500
$ws.Range("L1063").Value = 3500
$ws.Range("M1063").Value = 3500
$ws.Range("N1063").Value = "$/paquete 10 unidades (volumen en unidades)"
$ws.Range("O1063").Value = "Región de O'Higgins"
$ws.Range("P1063").Value = 350
$ws.Range("Q1063").Value = 10
$ws.Range("R1063").Value = "Hortaliza"

$ws.Range("A1064").Value = 5
$ws.Range("B1064").Value = "Macroferia Regional de Talca"
$ws.Range("C1064").Value = "Maule"
$ws.Range("D1064").Value = 45239
$ws.Range("E1064").Value = 7
$ws.Range("F1064").Value = 100112004
$ws.Range("G1064").Value = "Cebolla"
$ws.Range("H1064").Value = "Sin especificar"
$ws.Range("I1064").Value = "2a nueva(o)"
$ws.Range("J1064").Value = 10000
$ws.Range("K1064").Value = 2500
$ws.Range("L1064").Value = 2500
$ws.Range("M1064").Value = 2500
$ws.Range("N1064").Value = "$/paquete 10 unidades (volumen en unidades)"
$ws.Range("O1064").Value = "Región de O'Higgins"
$ws.Range("P1064").Value = 250
$ws.Range("Q1064").Value = 10
$ws.Range("R1064").Value = "Hortaliza"
